$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Set the B2:D9 block to 0 by default
$ws.Range("B2:D9").Value = 0

# Apply the specific non-zero values from the diff
$ws.Range("D2").Value = -0.6892262280447684
$ws.Range("D6").Value = 0.7001554245138226
